$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell in column B from "commentaire" to "nom".
# (Column A header "codeProgramme" stays the same.)
$ws.Range("B1").Value = "nom"

# Move/restore the active selection to B2 (was D18).
$ws.Range("B2").Select()
